$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet and name it
# "AutoWidthFalse" (the third case in the autoWidth fixture set).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "AutoWidthFalse"

# Same two-row/two-column data shape as the other AutoWidth* sheets,
# but with no explicit column widths set (the autoWidth=false default).
$ws.Range("A1").Value = "短"
$ws.Range("B1").Value = "長いヘッダーテキスト"
$ws.Range("A2").Value = "A"
$ws.Range("B2").Value = "データ"
